$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Add the two new backlog rows to the table. Using ListRows.Add() lets
#     the table (and therefore the sheet dimension / autofilter range)
#     grow the same way it does in the authored diff: A1:F119 -> A1:F121 ---

$newRow1 = $tbl.ListRows.Add()
$r1 = $newRow1.Range
$r1.Cells.Item(1,1).Value = "تعريف شيوه محاسبه و تعيين رده سني بازيكنان براساس استاندارد تعريف شده فدراسيون"
$r1.Cells.Item(1,2).Value = "سوم"
$r1.Cells.Item(1,3).Value = 0
$r1.Cells.Item(1,4).Value = 0
$r1.Cells.Item(1,5).Value = 0
$r1.Cells.Item(1,6).Value = 0

$newRow2 = $tbl.ListRows.Add()
$r2 = $newRow2.Range
$r2.Cells.Item(1,1).Value = "امكان به‌روز رساني دوره اي رده سني بازيكنان براساس استاندارد تعريف شده در نرم افزار"
$r2.Cells.Item(1,2).Value = "سوم"
$r2.Cells.Item(1,3).Value = 0
$r2.Cells.Item(1,4).Value = 0
$r2.Cells.Item(1,5).Value = 0
$r2.Cells.Item(1,6).Value = 0

# Column A keeps the wrap-text look used throughout the rest of the sheet;
# columns B:F keep the centered look used throughout the rest of the table.
$ws.Range("A120:A121").WrapText = $true
$ws.Range("B120:F121").HorizontalAlignment = -4108
$ws.Range("B120:F121").VerticalAlignment = -4108

# --- Move the active selection to mirror where the author was working ---
$ws.Range("A113").Select()
